# Generate Report for Handback
# Populates the "faac800b-c83b-4282-92e3-e43f2e2f61c8" handback row (row 8)
# on both locale sheets (zh-cn, de-de) now that a handback was processed:
#   - Latest Target File (I8): link to the handed-back markdown file
#   - Latest Handback File (J8): the generated xliff file name
#   - Latest Handback DateTime (K8): when the handback xliff was generated
#   - Error Detail (P8): the file was stale vs. the latest commit
# Also widens the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb2d4682d3d884f80ff5b063dc8e267566b83d49/e2e/faac800b-c83b-4282-92e3-e43f2e2f61c8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1b0a07343136bc10b60b8183740c83d7faf4356a/e2e/faac800b-c83b-4282-92e3-e43f2e2f61c8.md."

# ---------------- zh-cn sheet ----------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I8").Value = "faac800b-c83b-4282-92e3-e43f2e2f61c8.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I8"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/eb2d4682d3d884f80ff5b063dc8e267566b83d49/e2e/faac800b-c83b-4282-92e3-e43f2e2f61c8.md", "", "", "faac800b-c83b-4282-92e3-e43f2e2f61c8.md")

$wsZh.Range("J8").Value = "faac800b-c83b-4282-92e3-e43f2e2f61c8.bdaff405a56083d04b51715565b35ad23718e7f4.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-08-23 12:43:29"
$wsZh.Range("P8").Value = $errorDetail

# 39.17 "characters" renders as a stored sheet width of 40 (same calculation
# Excel already used for the other width=40 columns in this workbook).
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# ---------------- de-de sheet ----------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I8").Value = "faac800b-c83b-4282-92e3-e43f2e2f61c8.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I8"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/eb2d4682d3d884f80ff5b063dc8e267566b83d49/e2e/faac800b-c83b-4282-92e3-e43f2e2f61c8.md", "", "", "faac800b-c83b-4282-92e3-e43f2e2f61c8.md")

$wsDe.Range("J8").Value = "faac800b-c83b-4282-92e3-e43f2e2f61c8.bdaff405a56083d04b51715565b35ad23718e7f4.de-de.xlf"
$wsDe.Range("K8").Value = "2016-08-23 12:43:36"
$wsDe.Range("P8").Value = $errorDetail

$wsDe.Columns.Item(16).ColumnWidth = 39.17
